$d = $word.ActiveDocument

function Set-ParaFontSize10($para) {
    # Range.Text always includes the trailing paragraph-mark character
    # (\r), so a "visually empty" paragraph still reports Length 1.
    $txt = $para.Range.Text
    if ($txt.Length -gt 1) {
        # Paragraph already has run content: setting Font on the paragraph's
        # Range (which also spans its trailing paragraph mark) updates both
        # the run rPr's and the paragraph mark rPr (pPr/rPr) in one go.
        $para.Range.Font.Size = 10
        $para.Range.Font.SizeBi = 10
    } else {
        # Empty paragraph (no runs at all): a Range covering only the
        # paragraph mark cannot be reformatted directly. Insert a
        # placeholder character, format it (which stamps the mark's rPr
        # too), then delete the character again so the paragraph stays
        # empty but keeps the new rPr.
        $rng = $para.Range
        $rng.InsertBefore("X")
        $para.Range.Font.Size = 10
        $para.Range.Font.SizeBi = 10
        $delRng = $d.Range($para.Range.Start, $para.Range.Start + 1)
        $delRng.Delete()
    }
}

# Paragraphs 6-15 (1-based, Word numbering): the address block, "Dear ...",
# and the "In response to ..." approval intro -- all drop from 11pt (sz 22)
# to 10pt (sz 20).
for ($i = 6; $i -le 15; $i++) {
    Set-ParaFontSize10($d.Paragraphs.Item($i))
}

# The two bare "<w:p/>" paragraphs (no pPr at all) near the bookmark /
# signature area gain an explicit 10pt paragraph mark so they match their
# neighbours.
Set-ParaFontSize10($d.Paragraphs.Item(31))
Set-ParaFontSize10($d.Paragraphs.Item(34))
